$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 6447
$ws1.Range("F3").Value = 2589
$ws1.Range("F4").Value = 678
$ws1.Range("F5").Value = 102
$ws1.Range("F7").Value = 3180
$ws1.Range("F8").Value = 377
$ws1.Range("F11").Value = 8000
$ws1.Range("F12").Value = 409
$ws1.Range("F14").Value = 117
$ws1.Range("F18").Value = 494
$ws1.Range("F19").Value = 301
$ws1.Range("F20").Value = 9956
$ws1.Range("F22").Value = 277
$ws1.Range("F24").Value = 135
$ws1.Range("F25").Value = 376
$ws1.Range("F26").Value = 148
$ws1.Range("F30").Value = 123
$ws1.Range("F35").Value = 1500
$ws1.Range("F36").Value = 829
$ws1.Range("F37").Value = 4007
$ws1.Range("F38").Value = 244
$ws1.Range("F40").Value = 1718
$ws1.Range("F43").Value = 290
$ws1.Range("F44").Value = 175
$ws1.Range("F46").Value = 62
$ws1.Range("F47").Value = 53

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 6
$ws2.Range("F9").Value = 24
$ws2.Range("F13").Value = 20
$ws2.Range("F19").Value = 36
$ws2.Range("F20").Value = 16

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6447
$ws4.Range("F3").Value = 2589
$ws4.Range("F5").Value = 678
$ws4.Range("F6").Value = 102
$ws4.Range("F8").Value = 3180
$ws4.Range("F9").Value = 377
$ws4.Range("F13").Value = 8000
$ws4.Range("F14").Value = 409
$ws4.Range("F16").Value = 117
$ws4.Range("F20").Value = 301
$ws4.Range("F21").Value = 9956
$ws4.Range("F22").Value = 277
$ws4.Range("F24").Value = 135
$ws4.Range("F25").Value = 376
$ws4.Range("F26").Value = 148
$ws4.Range("F30").Value = 123
$ws4.Range("F33").Value = 1500
$ws4.Range("F34").Value = 829
$ws4.Range("F36").Value = 4007
$ws4.Range("F37").Value = 244
$ws4.Range("F39").Value = 1718
$ws4.Range("F43").Value = 290
$ws4.Range("F44").Value = 175
$ws4.Range("F46").Value = 62
$ws4.Range("F47").Value = 53
